$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "misc."

$ws.Range("K4").Value = "other_indexes"
$ws.Range("L4").Value = "commodity"
$ws.Range("K4:L4").Style = "Heading 3"

$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("L11").Value = "co2captured"
$ws.Range("K11").Value = "co2"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95

$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85
$ws.Range("K12").Value = "co2"
$ws.Range("L12").Value = "co2captured"

$ws.Range("D13").Select()

# Resize columns E and K to fit the new, wider content (mirrors Excel's
# "AutoFit column width" after the new cells were typed in).
$ws.Columns.Item(5).ColumnWidth = 9.46
$ws.Columns.Item(11).ColumnWidth = 11.1666666666666666
